# (v2.1.1.9051) fix VGS mo code
#
# The "microorganisms.codes" lookup table (code -> mo) on "Sheet 1" has a
# bogus row: code "VGS" incorrectly mapped to mo "B_VGCCC_SLMN". That whole
# row is spurious and must be removed entirely (not just blanked), so every
# row below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 0
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 1).Value() -eq "VGS") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows($targetRow).Delete()
}
